$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-6 from 45175 to 45183
$ws.Range("C2:C6").Value = 45183
